$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.242.08'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.875.33'
$ws.Range('E3').Value = '  -1.99%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.30'
$ws.Range('E5').Value = '  -1.77%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4692'
$ws.Range('E7').Value = '  -1.90%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2837'
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06605'
$ws.Range('E9').Value = '  -1.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.42'
$ws.Range('E10').Value = '  +8.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07779'
$ws.Range('E11').Value = '  +0.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '97.90'
$ws.Range('E12').Value = '  -4.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.878.14'
$ws.Range('E13').Value = '  -1.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.090'
$ws.Range('E14').Value = '  -2.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6739'
$ws.Range('E15').Value = '  +0.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '284.92'
$ws.Range('E16').Value = '  +7.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.257.89'
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9999'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.64'
$ws.Range('E19').Value = '  -0.44%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.398'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.122.95'
$ws.Range('E21').Value = '  -1.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.000007293'
$ws.Range('E22').Value = '  -2.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.182'
$ws.Range('E24').Value = '  -1.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.388'
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.60'
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.22'
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.988'
$ws.Range('E28').Value = '  -3.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.374'
$ws.Range('E29').Value = '  -0.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09679'
$ws.Range('E30').Value = '  -3.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.391'
$ws.Range('E31').Value = '  -4.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.471'
$ws.Range('E32').Value = '  -2.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.127'
$ws.Range('E33').Value = '  -2.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04694'
$ws.Range('E34').Value = '  -0.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7068'
$ws.Range('E35').Value = '  -2.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.093'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.719'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01879'
$ws.Range('E38').Value = '  -2.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.504'
$ws.Range('E39').Value = '  +3.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.525'
$ws.Range('E40').Value = '  -3.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.20'
$ws.Range('E41').Value = '  -3.72%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.966'
$ws.Range('E42').Value = '  -0.27%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8633'
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '103.15'
$ws.Range('E45').Value = '  -1.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4203'
$ws.Range('E46').Value = '  -1.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '986.05'
$ws.Range('E47').Value = '  +6.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.248'
$ws.Range('E48').Value = '  -2.32%  '
$ws.Range('E49').Value = '  +4.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.98'
$ws.Range('E50').Value = '  -2.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1151'
$ws.Range('E51').Value = '  -4.19%  '
